$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A3: 9999 -> -1
$ws.Range("A3").Value = -1

# A5: empty -> 9999 (keeps existing style s="5")
$ws.Range("A5").Value = 9999

# A6: empty -> shared string "null", styled like A4/A10 (style s="6")
$ws.Range("A6").Value = "null"
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A6").PasteSpecial(-4122) | Out-Null

# Update active selection from C6 to B4
$ws.Range("B4").Select() | Out-Null
